$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "4-3="
$t.Cell(1,2).Range.Text = "11+33="
$t.Cell(1,3).Range.Text = "92+1="
$t.Cell(1,4).Range.Text = "43-20="
$t.Cell(1,5).Range.Text = "8+76="

$t.Cell(2,1).Range.Text = "40-39="
$t.Cell(2,2).Range.Text = "12+66="
$t.Cell(2,3).Range.Text = "98-76="
$t.Cell(2,4).Range.Text = "84-26="
$t.Cell(2,5).Range.Text = "2+23="

$t.Cell(3,1).Range.Text = "87-7="
$t.Cell(3,2).Range.Text = "44+0="
$t.Cell(3,3).Range.Text = "0+10="
$t.Cell(3,4).Range.Text = "71+3="
$t.Cell(3,5).Range.Text = "38-7="

$t.Cell(4,1).Range.Text = "7+44="
$t.Cell(4,2).Range.Text = "78-25="
$t.Cell(4,3).Range.Text = "99-65="
$t.Cell(4,4).Range.Text = "76-63="
$t.Cell(4,5).Range.Text = "83-40="

$t.Cell(5,1).Range.Text = "0+1="
$t.Cell(5,2).Range.Text = "27+58="
$t.Cell(5,3).Range.Text = "73-27="
$t.Cell(5,4).Range.Text = "38-24="
$t.Cell(5,5).Range.Text = "81-47="

$t.Cell(6,1).Range.Text = "65-23="
$t.Cell(6,2).Range.Text = "88-44="
$t.Cell(6,3).Range.Text = "70-38="
$t.Cell(6,4).Range.Text = "67+19="
$t.Cell(6,5).Range.Text = "4+89="

$t.Cell(7,1).Range.Text = "65-28="
$t.Cell(7,2).Range.Text = "80+15="
$t.Cell(7,3).Range.Text = "32-28="
$t.Cell(7,4).Range.Text = "83-62="
$t.Cell(7,5).Range.Text = "56-7="

$t.Cell(8,1).Range.Text = "53-48="
$t.Cell(8,2).Range.Text = "15-6="
$t.Cell(8,3).Range.Text = "40-24="
$t.Cell(8,4).Range.Text = "74-42="
$t.Cell(8,5).Range.Text = "68+7="

$t.Cell(9,1).Range.Text = "82-58="
$t.Cell(9,2).Range.Text = "74+17="
$t.Cell(9,3).Range.Text = "2+5="
$t.Cell(9,4).Range.Text = "72-60="
$t.Cell(9,5).Range.Text = "71-64="

$t.Cell(10,1).Range.Text = "56+27="
$t.Cell(10,2).Range.Text = "5+53="
$t.Cell(10,3).Range.Text = "62-38="
$t.Cell(10,4).Range.Text = "43+49="
$t.Cell(10,5).Range.Text = "68+10="

$t.Cell(11,1).Range.Text = "79-76="
$t.Cell(11,2).Range.Text = "70+10="
$t.Cell(11,3).Range.Text = "27+9="
$t.Cell(11,4).Range.Text = "42-25="
$t.Cell(11,5).Range.Text = "8+18="

$t.Cell(12,1).Range.Text = "31-18="
$t.Cell(12,2).Range.Text = "14+82="
$t.Cell(12,3).Range.Text = "12+20="
$t.Cell(12,4).Range.Text = "46-22="
$t.Cell(12,5).Range.Text = "73+17="

$t.Cell(13,1).Range.Text = "13+81="
$t.Cell(13,2).Range.Text = "21+59="
$t.Cell(13,3).Range.Text = "72-17="
$t.Cell(13,4).Range.Text = "91-39="
$t.Cell(13,5).Range.Text = "60-10="

$t.Cell(14,1).Range.Text = "66-19="
$t.Cell(14,2).Range.Text = "6+21="
$t.Cell(14,3).Range.Text = "79-2="
$t.Cell(14,4).Range.Text = "6+54="
$t.Cell(14,5).Range.Text = "62-25="

$t.Cell(15,1).Range.Text = "99-68="
$t.Cell(15,2).Range.Text = "96-51="
$t.Cell(15,3).Range.Text = "83-2="
$t.Cell(15,4).Range.Text = "5+84="
$t.Cell(15,5).Range.Text = "49-18="

$t.Cell(16,1).Range.Text = "73-1="
$t.Cell(16,2).Range.Text = "67-16="
$t.Cell(16,3).Range.Text = "29+38="
$t.Cell(16,4).Range.Text = "98-89="
$t.Cell(16,5).Range.Text = "28+59="

$t.Cell(17,1).Range.Text = "89+10="
$t.Cell(17,2).Range.Text = "0+31="
$t.Cell(17,3).Range.Text = "34-9="
$t.Cell(17,4).Range.Text = "36+8="
$t.Cell(17,5).Range.Text = "86-10="

$t.Cell(18,1).Range.Text = "58-53="
$t.Cell(18,2).Range.Text = "74-51="
$t.Cell(18,3).Range.Text = "31+40="
$t.Cell(18,4).Range.Text = "25+20="
$t.Cell(18,5).Range.Text = "23+5="

$t.Cell(19,1).Range.Text = "33+3="
$t.Cell(19,2).Range.Text = "93+3="
$t.Cell(19,3).Range.Text = "92+0="
$t.Cell(19,4).Range.Text = "68+23="
$t.Cell(19,5).Range.Text = "67-19="

$t.Cell(20,1).Range.Text = "84-64="
$t.Cell(20,2).Range.Text = "42-15="
$t.Cell(20,3).Range.Text = "2+48="
$t.Cell(20,4).Range.Text = "57-26="
$t.Cell(20,5).Range.Text = "27+17="
